# menthol_CDCl3.xlsx: delete the now-unused H1_pureshift placeholder sheet
# and fill in the H1_1D sheet with the cached 1D-1H NMR peak-picking table
# (ppm / Integral / H's / Class / J's), matching the "started to add a hsqc
# CH only experiment, cached java output" commit.

$wb = $excel.ActiveWorkbook

# --- 1. drop the empty H1_pureshift sheet -----------------------------
$wb.Worksheets.Item("H1_pureshift").Delete() | Out-Null

# --- 2. populate H1_1D --------------------------------------------------
$ws = $wb.Worksheets.Item("H1_1D")

# Header row: B1:D1 (ppm/Intensity/Type) already exist; add Integral, H's,
# Class, J's in C1:F1 -- note Integral/H's slide into the old Intensity/
# Type slots, so set all four explicitly.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1:F1").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(1,2).Value = "ppm"
$ws.Cells.Item(1,3).Value = "Integral"
$ws.Cells.Item(1,4).Value = "H's"
$ws.Cells.Item(1,5).Value = "Class"
$ws.Cells.Item(1,6).Value = "J's"

# Column A (peak index) gets the same bold/boxed header style as the other
# sheets' index columns.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A11").PasteSpecial(-4122) | Out-Null

$peaks = @(
    @(1,  3.261234873010315,  1.016678020556872,  1, "td",    "4.35, 10.5"),
    @(2,  2.941423249842074,  0.9981795895695205, 1, "dd",    "1.88, 4.17"),
    @(3,  2.11101535336931,   1.00745855492204,    1, "heptd", "2.77, 6.98"),
    @(4,  1.858571924522433,  0.9351785302881365, 1, "dtd",   "2.08, 3.9, 12.3"),
    @(5,  1.522379543733942,  1.999049704789259,  2, "ddq",   "2.9, 13.3, 36.9"),
    @(6,  1.295270255624172,  0.7874049433813467, 1, "m",     $null),
    @(7,  1.004117802142161,  0.948374236143529,  1, "ddt",   "3.15, 10.1, 13.2"),
    @(8,  0.8490417558306967, 8.215744988315292,  8, "m",     $null),
    @(9,  0.7412999859978294, 1.038749497753711,  1, "ddd",   "3.56, 11.6, 12.9"),
    @(10, 0.6908204452587272, 3.053181934280298,  3, "d",     "7.22")
)

$row = 2
foreach ($p in $peaks) {
    $ws.Cells.Item($row, 1).Value = $p[0]
    $ws.Cells.Item($row, 2).Value = $p[1]
    $ws.Cells.Item($row, 3).Value = $p[2]
    $ws.Cells.Item($row, 4).Value = $p[3]
    $ws.Cells.Item($row, 5).Value = $p[4]
    if ($p[5] -ne $null) {
        $ws.Cells.Item($row, 6).Value = $p[5]
    }
    $row = $row + 1
}
